$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "TC_Functional_Smoke_002"
$ws.Range("B2").Value = "Pass"
$ws.Range("C2").Value = 45224.61497571259

# Row 3
$ws.Range("A3").Value = "TC_Functional_Smoke_003"
$ws.Range("B3").Value = "Pass"
$ws.Range("C3").Value = 45224.61656048317

# Row 4
$ws.Range("A4").Value = "TC_Functional_Smoke_008"
$ws.Range("B4").Value = "Pass"
$ws.Range("C4").Value = 45224.61818449155

# Row 5
$ws.Range("A5").Value = "TC_Functional_Smoke_009"
$ws.Range("B5").Value = "Pass"
$ws.Range("C5").Value = 45224.62033773593

# Row 6
$ws.Range("A6").Value = "TC_Functional_Smoke_010_1"
$ws.Range("B6").Value = "Pass"
$ws.Range("C6").Value = 45224.62299976052

# Row 7
$ws.Range("A7").Value = "obj.TC_Functional_Smoke_010_2"
$ws.Range("B7").Value = "Pass"
$ws.Range("C7").Value = 45224.62510186626

# Row 8
$ws.Range("A8").Value = "obj.TC_Functional_Smoke_018()"
$ws.Range("B8").Value = "Fail"
$ws.Range("C8").Value = 45224.63316037269

# Row 9
$ws.Range("A9").Value = "TC_Functional_Sanity_002_1"
$ws.Range("B9").Value = "Pass"
$ws.Range("C9").Value = 45224.63760603326

# Row 10
$ws.Range("A10").Value = "TC_Functional_Sanity_002_2"
$ws.Range("B10").Value = "Pass"
$ws.Range("C10").Value = 45224.64293216339

# Row 11
$ws.Range("A11").Value = "TC_Functional_Sanity_002_3"
$ws.Range("B11").Value = "Fail"
$ws.Range("C11").Value = 45224.64726967684

# Row 12 (new)
$ws.Range("A12").Value = "TC_Functional_Sanity_005"
$ws.Range("B12").Value = "Pass"
$ws.Range("C12").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C12").Value = 45224.64846922771
